$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# REPORTDATE: 2019-09-30 -> 2020-09-30
$ws.Range("H2").Value = "2020-09-30 00:00:00"

# BASIC_EPS
$ws.Range("I2").Value = 0.12

# TOTAL_OPERATE_INCOME
$ws.Range("K2").Value = 89942100.12

# PARENT_NETPROFIT
$ws.Range("L2").Value = 5357437.88

# YSTZ / SJLTZ / BPS / MGJYXJJE - were blank, now populated
$ws.Range("N2").Value = 27.5546035776
$ws.Range("O2").Value = -13.4346335287
$ws.Range("P2").Value = 0.630100861
$ws.Range("Q2").Value = 0.0091204535

# XSMLL
$ws.Range("R2").Value = 30.7711064374

# ISNEW: 0 -> 1 (kept as text, matching original inlineStr type)
$ws.Range("AB2").Value = "'1"

# QDATE
$ws.Range("AC2").Value = "2020Q3"

# DATATYPE
$ws.Range("AD2").Value = "2020年 三季报"

# DATAYEAR: 2019 -> 2020 (kept as text, matching original inlineStr type)
$ws.Range("AE2").Value = "'2020"
